$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix source probe: x-ray -> visible light
$ws.Range("J7").Value = "visible light"

# Fix image_key / monitor data dtype: uint32 -> int
$ws.Range("F3").Value = "int"
$ws.Range("F6").Value = "int"

# Update the active selection to match the saved view state
$ws.Range("D8").Select()
